# Insert 5 new publication rows at the top of the table (new rows 2-6),
# pushing the previously-first entry (CLOC / CVPR, etc.) down to row 7
# and all following rows down by 5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 blank rows before current row 2 (pushes existing data down)
$ws.Range("2:6").Insert()

# -----------------------------------------------------------------
# Populate the new string cells in the exact order the workbook's
# shared-string table records them (so new <si> entries line up).
# -----------------------------------------------------------------
$ws.Cells.Item(2,"B").Value = "Understanding Annotation Error Propagation and Learning an Adaptive Policy for Expert Intervention in Barrett's Video Segmentation"
$ws.Cells.Item(2,"C").Value = "Lokesha Rasanjalee, Dileepa Pitawela, Jin Tan, Rajvinder Singh, Tim Chen"
$ws.Cells.Item(2,"D").Value = "International Symposium on Biomedical Imaging (ISBI)"
$ws.Cells.Item(3,"B").Value = "OUGS: Active View Selection via Object-aware Uncertainty Estimation in 3DGS"
$ws.Cells.Item(3,"C").Value = "Haiyi Li , Qi Chen , Denis Kalkofe ,Tim Chen"
$ws.Cells.Item(3,"D").Value = "Eurographics"
$ws.Cells.Item(3,"H").Value = "Graphics,XR"
$ws.Cells.Item(3,"E").Value = "2026_EG_OUGS.jpg"
$ws.Cells.Item(2,"E").Value = "2026_ISBI_L2D.jpg"
$ws.Cells.Item(2,"F").Value = "2026_ISBI_L2D.pdf"
$ws.Cells.Item(4,"C").Value = "Carlos Tirado Cortes, Yiheng Chi, Juno Kim, Tim Chen"
$ws.Cells.Item(4,"D").Value = "Transaction on Visualisation and Computer Graphics (IEEE VR)"
$ws.Cells.Item(4,"E").Value = "2026_TVCG_Cybersickness.jpg"
$ws.Cells.Item(4,"F").Value = "2026_TVCG_Cybersickness.pdf"
$ws.Cells.Item(4,"B").Value = "Kinematic Sickness: Understanding Cybersickness Through Body Kinematics"
$ws.Cells.Item(5,"B").Value = "L2CU: Learning to Complement Unseen Users"
$ws.Cells.Item(5,"E").Value = "2025_Access_L2CU.jpg"
$ws.Cells.Item(5,"F").Value = "2025_Access_L2CU.pdf"
$ws.Cells.Item(6,"B").Value = "Learning To Defer To A Population With Limited Demonstrations"
$ws.Cells.Item(6,"D").Value = "DICTA"
$ws.Cells.Item(6,"E").Value = "2025_DICTA_L2D.jpg"
$ws.Cells.Item(6,"F").Value = "2025_DICTA_L2D.pdf"
$ws.Cells.Item(6,"C").Value = "Nilesh Ramgolam, Gustavo Carneiro, Tim Chen"

# -----------------------------------------------------------------
# Cells reusing strings already present in the shared-string table
# -----------------------------------------------------------------
$ws.Cells.Item(2,"H").Value = "AI, Medical"
$ws.Cells.Item(4,"H").Value = "XR"
$ws.Cells.Item(5,"C").Value = "Dileepa Pitawela, Gustavo Carneiro, Tim Chen"
$ws.Cells.Item(5,"D").Value = "IEEE Access"
$ws.Cells.Item(6,"H").Value = "AI"

# -----------------------------------------------------------------
# Year values (column A, numeric)
# -----------------------------------------------------------------
$ws.Cells.Item(2,"A").Value = 2026
$ws.Cells.Item(3,"A").Value = 2026
$ws.Cells.Item(4,"A").Value = 2026
$ws.Cells.Item(5,"A").Value = 2026
$ws.Cells.Item(6,"A").Value = 2025

# -----------------------------------------------------------------
# Row heights
# -----------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 48
$ws.Rows.Item(3).RowHeight = 32
$ws.Rows.Item(4).RowHeight = 32
$ws.Rows.Item(5).RowHeight = 16
$ws.Rows.Item(6).RowHeight = 16

# -----------------------------------------------------------------
# Wrap text (style index 1) on the B/C/D (title/authors/venue) columns
# -----------------------------------------------------------------
$ws.Range("B2:D6").WrapText = $true

# -----------------------------------------------------------------
# Restore the recorded active-cell selection
# -----------------------------------------------------------------
$ws.Range("C7").Select()
